# Import fund units2 — trims the capital-distributions sample sheet down
# to two generic "Distribution 1"/"Distribution 2" rows, clears the
# now-unused currency/commitment columns on those rows, and removes the
# third (CoInvest) sample row entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two sample distribution titles to the generic placeholders.
$ws.Range("C2").Value = "Distribution 1"
$ws.Range("C3").Value = "Distribution 2"

# Row 3's "Payments Paid" flips from No to Yes.
$ws.Range("I3").Value = "Yes"

# Distribution Basis / From Currency / To Currency / Exchange Rate / As Of
# are no longer populated for these rows - clear them (keeps any existing
# number formatting on the cells, e.g. the date format on column O).
$ws.Range("K2:O2").ClearContents()
$ws.Range("K3:O3").ClearContents()

# Drop the third sample row (the CoInvest / "Distribution from Sale : Nov" entry).
$ws.Rows("4:4").Delete()

# Matches the saved selection/view state in the target workbook.
$ws.Range("C4").Select()
